# Updated cryptos list on Tue Mar 19 16:00:27 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.939.64"
$ws.Range("E2").Value = "  -5.22%  "
$ws.Range("D3").Value = "3.308.18"
$ws.Range("E3").Value = "  -6.14%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "182.96"
$ws.Range("E5").Value = "  -10.00%  "
$ws.Range("D6").Value = "527.54"
$ws.Range("E6").Value = "  -4.61%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.304.18"
$ws.Range("E8").Value = "  -5.98%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "0.625"
$ws.Range("E10").Value = "  -5.22%  "
$ws.Range("D11").Value = "60.17"
$ws.Range("E11").Value = "  -2.05%  "
$ws.Range("E12").Value = "  -6.91%  "
$ws.Range("D13").Value = "0.0000261"
$ws.Range("E13").Value = "  -5.26%  "
$ws.Range("D14").Value = "9.19"
$ws.Range("E14").Value = "  -6.75%  "
$ws.Range("D15").Value = "3.832.71"
$ws.Range("E15").Value = "  -6.14%  "
$ws.Range("D16").Value = "0.118"
$ws.Range("E16").Value = "  -4.67%  "
$ws.Range("D17").Value = "3.308.98"
$ws.Range("E17").Value = "  -5.99%  "
$ws.Range("D18").Value = "17.73"
$ws.Range("E18").Value = "  -4.55%  "
$ws.Range("D19").Value = "63.946.95"
$ws.Range("E19").Value = "  -4.75%  "
$ws.Range("D20").Value = "11.12"
$ws.Range("E20").Value = "  -6.77%  "
$ws.Range("D21").Value = "0.967"
$ws.Range("E21").Value = "  -6.85%  "
$ws.Range("D22").Value = "373.63"
$ws.Range("E22").Value = "  -4.74%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").Value = "3.78"
$ws.Range("E23").Value = "  -6.66%  "
$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D24").Value = "11.29"
$ws.Range("E24").Value = "  -6.20%  "
$ws.Range("D25").Value = "81.12"
$ws.Range("E25").Value = "  -2.18%  "
$ws.Range("D26").Value = "3.96"
$ws.Range("E26").Value = "  +5.81%  "
$ws.Range("E27").Value = "  -1.48%  "
$ws.Range("D29").Value = "11.66"
$ws.Range("E29").Value = "  -3.57%  "
$ws.Range("D30").Value = "8.48"
$ws.Range("E30").Value = "  -5.23%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "29.02"
$ws.Range("E31").Value = "  -6.04%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "6.94"
$ws.Range("E32").Value = "  -5.38%  "
$ws.Range("D33").Value = "652.04"
$ws.Range("E33").Value = "  -5.20%  "
$ws.Range("D34").Value = "11.41"
$ws.Range("E34").Value = "  -3.53%  "
$ws.Range("E35").Value = "  -4.05%  "
$ws.Range("D36").Value = "59.36"
$ws.Range("E36").Value = "  -7.43%  "
$ws.Range("D37").Value = "0.403"
$ws.Range("E37").Value = "  -2.05%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "37.11"
$ws.Range("E39").Value = "  -6.41%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("E41").Value = "  -3.30%  "
$ws.Range("D42").Value = "2.931.37"
$ws.Range("E42").Value = "  -5.19%  "
$ws.Range("E43").Value = "  -4.32%  "
$ws.Range("D44").Value = "2.50"
$ws.Range("E44").Value = "  -2.75%  "
$ws.Range("E45").Value = "  -10.39%  "
$ws.Range("D46").Value = "2.93"
$ws.Range("E46").Value = "  +11.35%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("E48").Value = "  -2.30%  "
$ws.Range("E49").Value = "  -7.36%  "
$ws.Range("D50").Value = "0.127"
$ws.Range("E51").Value = "  +1.71%  "
